# Auto-generated Excel COM-interop script applying the Hades_Profits.xlsx diff
# Updates columns H-N (currentAveragePrice.. LeveProfitHQ) on specific rows
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3491.182
$ws.Range("I76").Value = 3455.4443
$ws.Range("J76").Value = 3652
$ws.Range("K76").Value = 3455.4443
$ws.Range("L76").Value = 3652
$ws.Range("M76").Value = -3140.4443
$ws.Range("N76").Value = -4282

$ws.Range("H79").Value = 3491.182
$ws.Range("I79").Value = 3455.4443
$ws.Range("J79").Value = 3652
$ws.Range("K79").Value = 3455.4443
$ws.Range("L79").Value = 3652
$ws.Range("M79").Value = -2363.4443
$ws.Range("N79").Value = -5836

$ws.Range("H112").Value = 37039084
$ws.Range("I112").Value = 200000770
$ws.Range("J112").Value = 2340.8635
$ws.Range("K112").Value = 600002310
$ws.Range("L112").Value = 7022.5905
$ws.Range("M112").Value = -600001202
$ws.Range("N112").Value = -9238.5905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1699.4634
$ws.Range("I2").Value = 1582.4333
$ws.Range("J2").Value = 2018.6364
$ws.Range("K2").Value = 1582.4333
$ws.Range("L2").Value = 2018.6364
$ws.Range("M2").Value = -1469.4333
$ws.Range("N2").Value = -2244.6364

$ws.Range("H32").Value = 23069.28
$ws.Range("I32").Value = 22657.633
$ws.Range("J32").Value = 24369.21
$ws.Range("K32").Value = 22657.633
$ws.Range("L32").Value = 24369.21
$ws.Range("M32").Value = -22370.633
$ws.Range("N32").Value = -24943.21

$ws.Range("H97").Value = 3290262.5
$ws.Range("I97").Value = 3677259.2
$ws.Range("J97").Value = 790
$ws.Range("K97").Value = 3677259.2
$ws.Range("L97").Value = 790
$ws.Range("M97").Value = -3676763.2
$ws.Range("N97").Value = -1782

$ws.Range("H116").Value = 1699.4634
$ws.Range("I116").Value = 1582.4333
$ws.Range("J116").Value = 2018.6364
$ws.Range("K116").Value = 1582.4333
$ws.Range("L116").Value = 2018.6364
$ws.Range("M116").Value = 711.5667000000001
$ws.Range("N116").Value = -6606.6364

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1699.4634
$ws.Range("I3").Value = 1582.4333
$ws.Range("J3").Value = 2018.6364
$ws.Range("K3").Value = 1582.4333
$ws.Range("L3").Value = 2018.6364
$ws.Range("M3").Value = -1468.4333
$ws.Range("N3").Value = -2246.6364

$ws.Range("H134").Value = 3573.2666
$ws.Range("I134").Value = 3806.7693
$ws.Range("J134").Value = 2055.5
$ws.Range("K134").Value = 11420.3079
$ws.Range("L134").Value = 6166.5
$ws.Range("M134").Value = -8885.3079
$ws.Range("N134").Value = -11236.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3756.2083
$ws.Range("I31").Value = 4158.3335
$ws.Range("J31").Value = 3622.1667
$ws.Range("K31").Value = 4158.3335
$ws.Range("L31").Value = 3622.1667
$ws.Range("M31").Value = -3863.3335
$ws.Range("N31").Value = -4212.1667

$ws.Range("H34").Value = 3756.2083
$ws.Range("I34").Value = 4158.3335
$ws.Range("J34").Value = 3622.1667
$ws.Range("K34").Value = 4158.3335
$ws.Range("L34").Value = 3622.1667
$ws.Range("M34").Value = -3956.3335
$ws.Range("N34").Value = -4026.1667

$ws.Range("H107").Value = 683.1429000000001
$ws.Range("I107").Value = 683.25
$ws.Range("J107").Value = 683
$ws.Range("K107").Value = 683.25
$ws.Range("L107").Value = 683
$ws.Range("M107").Value = 1236.75
$ws.Range("N107").Value = -4523

$ws.Range("H132").Value = 20050.309
$ws.Range("I132").Value = 1629.6487
$ws.Range("J132").Value = 57915
$ws.Range("K132").Value = 4888.9461
$ws.Range("L132").Value = 173745
$ws.Range("M132").Value = -2358.9461
$ws.Range("N132").Value = -178805

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 227.42857
$ws.Range("I23").Value = 201
$ws.Range("J23").Value = 231.83333
$ws.Range("K23").Value = 603
$ws.Range("L23").Value = 695.49999
$ws.Range("M23").Value = -368
$ws.Range("N23").Value = -1165.49999

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 0

$ws.Range("H107").Value = 660.95123
$ws.Range("I107").Value = 551.48
$ws.Range("J107").Value = 832
$ws.Range("K107").Value = 1654.44
$ws.Range("L107").Value = 2496
$ws.Range("M107").Value = 265.5599999999999
$ws.Range("N107").Value = -6336

$ws.Range("H113").Value = 710.94116
$ws.Range("I113").Value = 600
$ws.Range("J113").Value = 757.1667
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 2271.5001
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -6611.5001

$ws.Range("H122").Value = 817.75
$ws.Range("I122").Value = 328.36365
$ws.Range("J122").Value = 1074.0952
$ws.Range("K122").Value = 2955.27285
$ws.Range("L122").Value = 9666.8568
$ws.Range("M122").Value = -505.2728500000003
$ws.Range("N122").Value = -14566.8568

$ws.Range("H131").Value = 885.8570999999999
$ws.Range("I131").Value = 484.83334
$ws.Range("J131").Value = 1186.625
$ws.Range("K131").Value = 1454.50002
$ws.Range("L131").Value = 3559.875
$ws.Range("M131").Value = 3585.49998
$ws.Range("N131").Value = -13639.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 36811
$ws.Range("I51").Value = 20296
$ws.Range("J51").Value = 53326
$ws.Range("K51").Value = 20296
$ws.Range("L51").Value = 53326
$ws.Range("M51").Value = -19787
$ws.Range("N51").Value = -54344

$ws.Range("H70").Value = 82659.234
$ws.Range("I70").Value = 337866.66
$ws.Range("J70").Value = 6097
$ws.Range("K70").Value = 337866.66
$ws.Range("L70").Value = 6097
$ws.Range("M70").Value = -337596.66
$ws.Range("N70").Value = -6637

$ws.Range("H73").Value = 82659.234
$ws.Range("I73").Value = 337866.66
$ws.Range("J73").Value = 6097
$ws.Range("K73").Value = 337866.66
$ws.Range("L73").Value = 6097
$ws.Range("M73").Value = -336930.66
$ws.Range("N73").Value = -7969

$ws.Range("H122").Value = 4163
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4163
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = 12489
$ws.Range("N122").Value = -17389

$ws.Range("H126").Value = 1542.4
$ws.Range("I126").Value = 1542.4
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4627.200000000001
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -2157.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 170499.67
$ws.Range("I132").Value = 2500
$ws.Range("J132").Value = 254499.5
$ws.Range("K132").Value = 7500
$ws.Range("L132").Value = 763498.5
$ws.Range("M132").Value = -4970
$ws.Range("N132").Value = -768558.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 40016.75
$ws.Range("I132").Value = 26242.025
$ws.Range("J132").Value = 85932.5
$ws.Range("K132").Value = 78726.07500000001
$ws.Range("L132").Value = 257797.5
$ws.Range("M132").Value = -76196.07500000001
$ws.Range("N132").Value = -262857.5
